$d = $word.ActiveDocument

# Build the WordprocessingML fragment for the four new paragraphs that need
# to be inserted at the very top of the document body:
#   1) Bold title-style paragraph: "Acknowledgement of Generative AI"
#   2) Body paragraph explaining generative AI usage
#   3) Body paragraph pointing to the chat links below
#   4) Empty bold/title-style paragraph (spacer)
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$titleRPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-GB"/></w:rPr>'
$bodyRPr  = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-GB"/></w:rPr>'

$p1 = '<w:p ' + $ns + '><w:pPr>' + $titleRPr + '</w:pPr><w:r>' + $titleRPr + '<w:t>Acknowledgement of Generative AI</w:t></w:r></w:p>'
$p2 = '<w:p ' + $ns + '><w:pPr>' + $bodyRPr + '</w:pPr><w:r>' + $bodyRPr + '<w:t>Throughout this project we used generative AI (ChatGPT) to help formulate the idea and most importantly create a proof-of-concept version of our platform Data Centre Frontier.</w:t></w:r></w:p>'
$p3 = '<w:p ' + $ns + '><w:pPr>' + $bodyRPr + '</w:pPr><w:r>' + $bodyRPr + '<w:t>Below you will find the links to the chats used by both of us.</w:t></w:r></w:p>'
$p4 = '<w:p ' + $ns + '><w:pPr>' + $titleRPr + '</w:pPr></w:p>'

$xml = $p1 + $p2 + $p3 + $p4

# Insert the fragment immediately before the current start of the document.
$r = $d.Range(0, 0)
$r.InsertXML($xml)
